$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data rows
$data = @(
  @(45620.99999999999, 24, 23.9999634925206,  23.99996349702972),
  @(45634.99999999999, 20, 19.99995548871507, 19.9999554934065),
  @(45641.99999999999, 18, 17.99995144260024, 17.99995154341689),
  @(45648.99999999999, 16, 15.99994732710003, 15.99994766659238),
  @(45655.99999999999, 14, 13.99994315590351, 13.99994383249053),
  @(45662.99999999999, 12, 11.99993891941044, 11.99994007663973),
  @(45669.99999999999, 10, 9.999934650496472, 9.999936330449271),
  @(45676.99999999999, 8,  7.999930376401295, 7.999932645846862),
  @(45683.99999999999, 6,  5.99992597908883,  5.999929017743712),
  @(45690.99999999999, 4,  3.999921636243711, 3.99992534517589)
)

$r = 2
foreach ($row in $data) {
  $wsForecast.Cells.Item($r, 1).Value = $row[0]
  $wsForecast.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
  $wsForecast.Cells.Item($r, 2).Value = $row[1]
  $wsForecast.Cells.Item($r, 3).Value = $row[2]
  $wsForecast.Cells.Item($r, 4).Value = $row[3]
  $r++
}

[void]$wsForecast.Range("A1").Select()
